$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.849.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.400.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.08%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.40%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'561.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.84%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'141.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.32%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.93%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.406.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.11%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.44%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.28%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.46%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.817.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.89%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'60.454.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.68%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.440.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +6.47%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.27%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'323.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.06%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.51%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.71%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'64.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.31%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'573.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.76%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E30").Value = "'  +0.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.10%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.64%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.05%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.59%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'153.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.83%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.35%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.16%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +8.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'41.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'142.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.33%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.48%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.79%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'19.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.95%  "
$ws.Range("E51").Style = "Normal"
